$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.353.50"
$ws.Range("E2").Value = "  -4.10%  "
$ws.Range("D3").Value = "1.764.68"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'304.61"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("D7").Value = "'0.4290"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").Value = "'0.3611"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "'0.07045"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'0.8337"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").Value = "1.740.10"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'6.419"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.233"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").Value = "'0.06786"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'79.01"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "'0.000008593"
$ws.Range("E18").Value = "  -3.48%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'14.96"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "26.365.17"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").Value = "'5.005"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "'11.08"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "1.986.35"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Value = "'152.63"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "'1.855"
$ws.Range("E26").Value = "  -6.44%  "
$ws.Range("D27").Value = "'18.10"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "'5.059"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "'114.34"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "'1.694"
$ws.Range("E30").Value = "  -5.38%  "
$ws.Range("D31").Value = "'0.08910"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "'0.7246"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "'4.322"
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("D34").Value = "'1.104"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "'2.749"
$ws.Range("E35").Value = "  -7.71%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "'1.071"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "'0.05107"
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("D39").Value = "'0.01887"
$ws.Range("D40").Value = "'0.4907"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").Value = "'0.1598"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  -9.56%  "
$ws.Range("D43").Value = "'6.202"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("D44").Value = "'8.021"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").Value = "'104.97"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'10.05"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("D48").Value = "'0.06189"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").Value = "'0.4471"
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("D50").Value = "'1.573"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "'1.721"
$ws.Range("E51").Value = "  -0.05%  "
